# Insert a new data row at row 174 (shifting existing rows 174-274 down to 175-275)
# and populate it with the new weekly price record, matching the rest of the
# row's fields (Mercado, Categoria, etc.) that stay identical to the row that
# used to occupy position 174.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(174).Insert()

$ws.Range("A174").Value = 9
$ws.Range("B174").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C174").Value = "Metropolitana"
$ws.Range("D174").Value = 44719
$ws.Range("E174").Value = 13
$ws.Range("F174").Value = 300000001
$ws.Range("G174").Value = "Rabanito"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 7000
$ws.Range("K174").Value = 3000
$ws.Range("L174").Value = 3000
$ws.Range("M174").Value = 3000
$ws.Range("N174").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O174").Value = "Provincia de Chacabuco"
$ws.Range("P174").Value = 30
$ws.Range("Q174").Value = 100
$ws.Range("R174").Value = "Hortaliza"
